# Insert two new weekly records for "Crespo record" (and the row that used to
# sit at row 51 onward all shift down by two rows to make room), then populate
# the two freshly inserted rows (51 and 52) with their own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows above the current row 51; this pushes the existing rows
# 51..126 down to 53..128 (and copies the formatting - e.g. the date style
# on column D - from the row being pushed down, same as a real Excel insert).
$ws.Range("A51:A52").EntireRow.Insert()

# New row 51: Crespo record / Primera
$ws.Cells.Item(51, 1).Value  = 5
$ws.Cells.Item(51, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(51, 3).Value  = "Maule"
$ws.Cells.Item(51, 4).Value  = 44413
$ws.Cells.Item(51, 5).Value  = 7
$ws.Cells.Item(51, 6).Value  = 100112006
$ws.Cells.Item(51, 7).Value  = "Repollo"
$ws.Cells.Item(51, 8).Value  = "Crespo record"
$ws.Cells.Item(51, 9).Value  = "Primera"
$ws.Cells.Item(51, 10).Value = 3000
$ws.Cells.Item(51, 11).Value = 500
$ws.Cells.Item(51, 12).Value = 500
$ws.Cells.Item(51, 13).Value = 500
$ws.Cells.Item(51, 14).Value = "$/unidad"
$ws.Cells.Item(51, 15).Value = "Región del Maule"
$ws.Cells.Item(51, 16).Value = 500
$ws.Cells.Item(51, 17).Value = 1
$ws.Cells.Item(51, 18).Value = "Hortaliza"

# New row 52: Crespo record / Segunda
$ws.Cells.Item(52, 1).Value  = 5
$ws.Cells.Item(52, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(52, 3).Value  = "Maule"
$ws.Cells.Item(52, 4).Value  = 44413
$ws.Cells.Item(52, 5).Value  = 7
$ws.Cells.Item(52, 6).Value  = 100112006
$ws.Cells.Item(52, 7).Value  = "Repollo"
$ws.Cells.Item(52, 8).Value  = "Crespo record"
$ws.Cells.Item(52, 9).Value  = "Segunda"
$ws.Cells.Item(52, 10).Value = 3000
$ws.Cells.Item(52, 11).Value = 350
$ws.Cells.Item(52, 12).Value = 350
$ws.Cells.Item(52, 13).Value = 350
$ws.Cells.Item(52, 14).Value = "$/unidad"
$ws.Cells.Item(52, 15).Value = "Región del Maule"
$ws.Cells.Item(52, 16).Value = 350
$ws.Cells.Item(52, 17).Value = 1
$ws.Cells.Item(52, 18).Value = "Hortaliza"
